$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update 최종점수 (K column) values
$ws.Range("K2").Value = 58.3
$ws.Range("K3").Value = 52.7

# Update MACRO_SCORE (N column) values
$ws.Range("N2").Value = 51.15965480231979
$ws.Range("N3").Value = 51.15965480231979
